$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "1.176") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.969.87"
$ws.Range("E2").Value = "  +5.98%  "

$ws.Range("D3").Value = "1.711.97"
$ws.Range("E3").Value = "  +3.66%  "

$ws.Range("D4").Value = "0.9974"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").Value = "331.66"
$ws.Range("E5").Value = "  +5.89%  "

$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "0.3694"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("D8").Value = "48.65"
$ws.Range("E8").Value = "  +4.42%  "

$ws.Range("D9").Value = "0.3316"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("D10").Value = "1.176"
$ws.Range("E10").Value = "  +4.53%  "

$ws.Range("D11").Value = "0.07418"
$ws.Range("E11").Value = "  +5.53%  "

$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "6.231"
$ws.Range("E13").Value = "  +4.46%  "

$ws.Range("D14").Value = "20.08"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").Value = "6.899"
$ws.Range("E15").Value = "  +4.29%  "

$ws.Range("D16").Value = "1.700.62"
$ws.Range("E16").Value = "  +2.88%  "

$ws.Range("D17").Value = "0.00001073"
$ws.Range("E17").Value = "  +3.04%  "

$ws.Range("D18").Value = "0.06654"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").Value = "81.46"
$ws.Range("E19").Value = "  +3.54%  "

$ws.Range("D20").Value = "0.9985"

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "16.24"
$ws.Range("E21").Value = "  +3.67%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.077"
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").Value = "13.04"
$ws.Range("E23").Value = "  +4.00%  "

$ws.Range("D24").Value = "25.824.76"
$ws.Range("E24").Value = "  +5.47%  "

$ws.Range("D25").Value = "2.471"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").Value = "2.489"
$ws.Range("E26").Value = "  +7.14%  "

$ws.Range("D27").Value = "149.82"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").Value = "19.21"
$ws.Range("E28").Value = "  +3.40%  "

$ws.Range("D29").Value = "1.300"
$ws.Range("E29").Value = "  +9.73%  "

$ws.Range("D30").Value = "1.892.48"
$ws.Range("E30").Value = "  +3.09%  "

$ws.Range("D31").Value = "129.04"
$ws.Range("E31").Value = "  +3.88%  "

$ws.Range("D32").Value = "4.104"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").Value = "5.976"
$ws.Range("E33").Value = "  +4.45%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "0.08527"
$ws.Range("E34").Value = "  +0.93%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.720"
$ws.Range("E35").Value = "  +3.67%  "

$ws.Range("D36").Value = "13.00"
$ws.Range("E36").Value = "  +6.81%  "

$ws.Range("D37").Value = "5.350"
$ws.Range("E37").Value = "  +2.83%  "

$ws.Range("D38").Value = "1.273"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").Value = "0.06217"
$ws.Range("E39").Value = "  +3.31%  "

$ws.Range("D40").Value = "8.543"
$ws.Range("E40").Value = "  +5.32%  "

$ws.Range("D41").Value = "0.2125"
$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("D42").Value = "0.02268"
$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("D43").Value = "14.63"
$ws.Range("E43").Value = "  +15.32%  "

$ws.Range("D44").Value = "0.6131"
$ws.Range("E44").Value = "  +4.04%  "

$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D47").Value = "0.5858"
$ws.Range("E47").Value = "  +4.32%  "

$ws.Range("D48").Value = "127.09"
$ws.Range("E48").Value = "  +2.99%  "

$ws.Range("D49").Value = "2.007"
$ws.Range("E49").Value = "  +3.02%  "

$ws.Range("D50").Value = "0.07238"
$ws.Range("E50").Value = "  +4.68%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "76.75"
$ws.Range("E51").Value = "  +3.16%  "
